# Update tab names in all BOMs, fix bi-color LED naming.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "iVCA" to "BOM"
$ws.Name = "BOM"

# Fix the LED package naming: "LED 3mm Flat Bicolor" -> "LED 3mm Dome Bicolor"
$ws.Range("C19").Value = "LED 3mm Dome Bicolor"

# Leave the selection on the cell that was edited, matching the saved view state
$ws.Range("C19").Select()
